# Add "Execution Time (ms)" and "Memory Usage (B)" metrics columns
# for each model in the test-set metrics table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Reuse the same header formatting (bold, centered, bordered) as the
# existing headers (e.g. D1) by copying just the formats over.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows -----------------------------------------------------------
# Row 2 - Linear
$ws.Range("E2").Value = 7.445600000210106
$ws.Range("F2").Value = 53248

# Row 3 - Decision Tree
$ws.Range("E3").Value = 7.472499972209334
$ws.Range("F3").Value = 53248

# Row 4 - Random Forest
$ws.Range("E4").Value = 20.87900001788512
$ws.Range("F4").Value = 0

# Row 5 - Lasso
$ws.Range("E5").Value = 5.950500024482608
$ws.Range("F5").Value = 4096

# Row 6 - Optimized Equation
$ws.Range("E6").Value = 4.208099999232218
$ws.Range("F6").Value = 4096
